$wb = $excel.ActiveWorkbook

# Row 2 values for columns A..O, one array per sheet, in sheet order
# (sheet 1 = "2025", sheet 2 = "2030", sheet 3 = "2035",
#  sheet 4 = "2040", sheet 5 = "2045", sheet 6 = "2050")
$sheetData = @(
    @(0, 1037.265132737054, 0, 0, 28926.05393052954, 0, 8095.925712661834, 0, 16171.06685703679, 0, 0, 48492.22142001599, 10595.37713982, 7071.74531360843, 6993.890772562212),
    @(0, 4157.588990853394, 0, 0, 45991.90904307188, 0, 8095.925712661834, 0, 37079.12819938764, 0, 0, 54844.03303316472, 17449.04999683176, 9024.733389685653, 9724.258249348202),
    @(2754.31755456332, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13034.3101291405, 12860.17168993684),
    @(2754.31755456332, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13151.8694977663, 12860.17168993684),
    @(5713.151062849596, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13601.08685191924, 14937.1305943757),
    @(5713.151062849596, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13601.08685191924, 14937.1305943757)
)

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $values = $sheetData[$s - 1]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item(2, $col).Value = $values[$i]
    }
}
